$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Scalar")
$ws.Range("B5").Value = 16679742.792528583
$ws.Range("B12").Value = 10811543.460479999
$ws.Range("C22").Value = 11495321.064
$ws.Range("C23").Value = 16679742.792528583
$ws.Range("C25").Value = 10811543.460479999
$ws.Range("B32").Value = -30439
$ws.Range("C32").Value = -30439
$ws.Range("D32").Value = -30439
$ws.Range("C33").Value = 10811543.460479999
$ws.Range("C34").Value = 10811543.460479999
$ws.Range("B36").Value = 11495321.064
$ws.Range("C36").Value = 11495321.064
$ws.Range("D36").Value = 11495321.064

$ws = $wb.Worksheets.Item("EnergyRate")
$ws.Range("B4").Value = 47.189
$ws.Range("B5").Value = 65.218

$ws = $wb.Worksheets.Item("EQ8_SteadyEnergy")
$ws.Range("E4").Value = 377.512
$ws.Range("E5").Value = 377.512

$ws = $wb.Worksheets.Item("EQ9_EnergyGen")
$ws.Range("E4").Value = 1085.347
$ws.Range("E5").Value = 1500.014

$ws = $wb.Worksheets.Item("evap")
$ws.Range("B4").Value = 986
$ws.Range("B5").Value = 986
$ws.Range("B6").Value = 985
$ws.Range("B7").Value = 984
$ws.Range("B8").Value = 983
$ws.Range("B9").Value = 983
$ws.Range("B10").Value = 983
$ws.Range("B11").Value = 984
$ws.Range("B12").Value = 984
$ws.Range("B13").Value = 984
$ws.Range("B14").Value = 985
$ws.Range("B15").Value = 985
$ws.Range("B16").Value = 985
$ws.Range("B17").Value = 984
$ws.Range("B18").Value = 984
$ws.Range("B19").Value = 983
$ws.Range("B20").Value = 983
$ws.Range("B21").Value = 983
$ws.Range("B22").Value = 982
$ws.Range("B23").Value = 982
$ws.Range("B24").Value = 981
$ws.Range("B25").Value = 981
$ws.Range("B26").Value = 980
$ws.Range("B27").Value = 980
$ws.Range("B28").Value = 979
$ws.Range("B29").Value = 978
$ws.Range("B30").Value = 978
$ws.Range("B31").Value = 977
$ws.Range("B32").Value = 976
$ws.Range("B33").Value = 976
$ws.Range("B34").Value = 975

$ws = $wb.Worksheets.Item("Inflow")
$ws.Range("B4").Value = 4614
$ws.Range("B5").Value = 5346
$ws.Range("B6").Value = 4425
$ws.Range("B7").Value = 4396
$ws.Range("B8").Value = 7281
$ws.Range("B9").Value = 7230
$ws.Range("B10").Value = 13082
$ws.Range("B11").Value = 14844
$ws.Range("B12").Value = 11161
$ws.Range("B13").Value = 13579
$ws.Range("B14").Value = 16356
$ws.Range("B15").Value = 7879
$ws.Range("B16").Value = 9725
$ws.Range("B17").Value = 8960
$ws.Range("B18").Value = 4362
$ws.Range("B19").Value = 5885
$ws.Range("B20").Value = 9247
$ws.Range("B21").Value = 9919
$ws.Range("B22").Value = 5928
$ws.Range("B23").Value = 8213
$ws.Range("B24").Value = 8163
$ws.Range("B25").Value = 9233
$ws.Range("B26").Value = 3347
$ws.Range("B27").Value = 6373
$ws.Range("B28").Value = 6926
$ws.Range("B29").Value = 5578
$ws.Range("B30").Value = 5377
$ws.Range("B31").Value = 6117
$ws.Range("B32").Value = 5425
$ws.Range("B33").Value = 5949
$ws.Range("B34").Value = 5572

$ws = $wb.Worksheets.Item("weekendRate")
$ws.Range("B4").Value = 47.189
$ws.Range("B5").Value = 47.189
